# Update the "welkomst" planning sheet: a handful of tasks that were still
# outstanding ("Nee") are now marked as finished ("Ja") with a completion
# date, and the remembered active cell/selection moves accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Rows whose "Afgemaakt" status flips from "Nee" to "Ja" with a finish date
$rows = @(10, 14, 19, 21, 22)

# Reuse the formatting already used for other completed rows (e.g. row 7)
# so the date column keeps its date number format and the status column
# keeps the "done" (green) fill instead of the "not done" (red) fill.
$doneDateFormat = $ws.Range("H7").NumberFormat
$doneFillColor = $ws.Range("I7").Interior.Color

foreach ($r in $rows) {
    $hCell = $ws.Range("H$r")
    $iCell = $ws.Range("I$r")

    $hCell.NumberFormat = $doneDateFormat
    $hCell.Value = 43391

    $iCell.Value = "Ja"
    $iCell.Interior.Color = $doneFillColor
}

# Move the active selection/cell as recorded in the saved view state
$ws.Range("G11").Select()
